$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 52

# A52: literal text "2026/01/01" (must NOT be auto-converted to a date serial).
# Enter it as a text formula, then convert the formula to a static value via
# a values-only paste so the stored cell keeps t="s" (string) with no formula
# residue and no NumberFormat change.
$a = $ws.Cells.Item($row, 1)
$a.Formula = "=""2026/01/01"""
$a.Copy()
$a.PasteSpecial(-4163)  # xlPasteValues

$ws.Cells.Item($row, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($row, 3).Value = 1123

# Copy the formatting (center/center alignment = style used by the rest of
# the data rows) from the previous row onto the new one in a single shot so
# no extra/unused style entries are minted.
$srcFmt = $ws.Range("A51:C51")
$srcFmt.Copy()
$dstFmt = $ws.Range("A$($row):C$($row)")
$dstFmt.PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
